$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fresh-Rotten")

# New "Aliment" column header and existing rows 2-4 tagged as "all"
$ws.Range("J1").Value = "Aliment"
$ws.Range("J2").Value = "all"
$ws.Range("J3").Value = "all"
$ws.Range("J4").Value = "all"

# New data rows 5-7 (apple)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.59337349397590367
$ws.Range("H5").Value = 344
$ws.Range("I5").Value = $false
$ws.Range("J5").Value = "apple"

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 0.91767068273092367
$ws.Range("H6").Value = 802
$ws.Range("I6").Value = $false
$ws.Range("J6").Value = "apple"

$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 0.93574297188755018
$ws.Range("H7").Value = 1165
$ws.Range("I7").Value = $false
$ws.Range("J7").Value = "apple"

# Row 8 - added without the "Aliment" tag
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = 0.995609220636663
$ws.Range("H8").Value = 1195
$ws.Range("I8").Value = $false

# Touch K8 so the sheet's used range extends to column K without changing styles
$ws.Range("K8").Style = "Normal"

$ws.Range("H8").Select()
